# Add the "SWC Clan Wars" channel to the CoC channel list (Sheet1, row 32).
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Fill in the new row's cells. Column B (channel id) is written before column A
# (channel name) so the shared-string table picks up the same ordering Excel
# produced for this edit.
$ws.Range("B32").Value = "UC_lPhq9a6rG76HXRFCl0zAw"
$ws.Range("A32").Value = "SWC Clan Wars"
$ws.Range("C32").Value = "한달후 확인"

# Column D holds a "recheck" date for this entry.
$ws.Range("D32").Value = 42124
$ws.Range("D32").NumberFormat = "m/d/yyyy"

# Column D needs to be wide enough to show the date.
$ws.Columns.Item(4).ColumnWidth = 10.4

# Update the view so the new row is in frame / selected, like after the edit.
$ws.Range("D33").Select() | Out-Null
